$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Food-Beverages: zoom 130 -> 220, active cell -> C3
# ---------------------------------------------------------------------------
$wsFB = $wb.Worksheets.Item("Food-Beverages")
$wsFB.Activate()
$excel.ActiveWindow.Zoom = 220
$wsFB.Range("C3").Select() | Out-Null

# ---------------------------------------------------------------------------
# Apparel: active cell -> C9 (tabSelected naturally moves off this sheet
# once a later sheet is activated below)
# ---------------------------------------------------------------------------
$wsAP = $wb.Worksheets.Item("Apparel")
$wsAP.Activate()
$wsAP.Range("C9").Select() | Out-Null

# ---------------------------------------------------------------------------
# Retail: active cell -> J30. The external-link-driven model refresh
# (GME.xlsx) updated cached prices; since that external workbook isn't
# reachable from this sandbox, push the refreshed numbers in directly.
# F30 and H30 pull straight from the external link, so they become literal
# values; E30 (=D30*H30) and G30 (=E30-F30) are purely local formulas, so
# they recalculate correctly on their own once F30/H30 are updated.
# ---------------------------------------------------------------------------
$wsRT = $wb.Worksheets.Item("Retail")
$wsRT.Activate()
$wsRT.Range("F30").Value = 4169.7000000000007
$wsRT.Range("H30").Value = 426.509592
$wsRT.Range("J30").Select() | Out-Null

# ---------------------------------------------------------------------------
# Leisure: active cell -> C3
# ---------------------------------------------------------------------------
$wsLE = $wb.Worksheets.Item("Leisure")
$wsLE.Activate()
$wsLE.Range("C3").Select() | Out-Null

# ---------------------------------------------------------------------------
# NonDurable: widen Name/Ticker columns, add right-aligned Price entries for
# a handful of names, and move the active cell -> B16.
# ---------------------------------------------------------------------------
$wsND = $wb.Worksheets.Item("NonDurable")
$wsND.Activate()

$wsND.Columns("B").ColumnWidth = 21.85
$wsND.Columns("C").ColumnWidth = 7.85

$wsND.Range("D2").HorizontalAlignment = -4152
$wsND.Range("E2").HorizontalAlignment = -4152

$wsND.Range("D3").Value = 172.31
$wsND.Range("D3").HorizontalAlignment = -4152

$wsND.Range("D5").Value = 119.72
$wsND.Range("D5").HorizontalAlignment = -4152

$wsND.Range("D8").Value = 85.5
$wsND.Range("D8").NumberFormat = "#,##0.00"
$wsND.Range("D8").HorizontalAlignment = -4152

$wsND.Range("D9").Value = 50.14
$wsND.Range("D9").HorizontalAlignment = -4152

$wsND.Range("D10").Value = 102.37
$wsND.Range("D10").HorizontalAlignment = -4152

$wsND.Range("D15").Value = 102.04
$wsND.Range("D15").HorizontalAlignment = -4152

$wsND.Range("B16").Select() | Out-Null

# ---------------------------------------------------------------------------
# Retail is the sheet left active/selected in the saved workbook.
# ---------------------------------------------------------------------------
$wsRT.Activate()
